$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("design","Test Prop 1","Horace's Home","2016-12-05","2016-12-08"),
    @("design","Test Prop 2","Horace's Home","2016-12-13","2016-12-16"),
    @("first build","Test Prop 1","Horace's Home","2016-12-19","2017-01-23"),
    @("first build","Test Prop 2","Horace's Home","2016-12-23","2017-01-23"),
    @("design","Test Prop 2","Horace's Home","2017-01-09","2017-01-12"),
    @("first build","Test Prop 3","Horace's Home","2017-01-23","2017-02-10"),
    @("first build","Test Prop 4","Not available","2017-01-23","2017-02-07"),
    @("first build","Test Prop 5","Horace's Home","2017-01-26","2017-02-08"),
    @("duplicate build","Test Prop 1","Horace's Home","2017-03-20","2017-03-24"),
    @("duplicate build","Test Prop 2","Horace's Home 2","2017-03-20","2017-03-24")
)

# Ensure the date columns are treated as text cells so the date-like strings
# are stored as literal text, not converted to Excel date serial numbers.
$ws.Range("D2:E11").NumberFormat = "@"

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r++
}
